$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cells for the new latitude/longitude columns
$ws.Range("E1").Value = "latitude"
$ws.Range("F1").Value = "longitude"

# Per-shop latitude / longitude values (rows 2-31)
$ws.Range("E2").Value = 40.440624999999997
$ws.Range("F2").Value = -79.995885999999999
$ws.Range("E3").Value = 42.443961000000002
$ws.Range("F3").Value = -76.501880999999997
$ws.Range("E4").Value = 39.961176000000002
$ws.Range("F4").Value = -82.998794000000004
$ws.Range("E5").Value = 42.331426999999998
$ws.Range("F5").Value = -83.045754000000002
$ws.Range("E6").Value = 39.103118000000002
$ws.Range("F6").Value = -84.512020000000007
$ws.Range("E7").Value = 38.252665
$ws.Range("F7").Value = -85.758455999999995
$ws.Range("E8").Value = 36.162663999999999
$ws.Range("F8").Value = -86.781602000000007
$ws.Range("E9").Value = 39.739235999999998
$ws.Range("F9").Value = -104.990251
$ws.Range("E10").Value = 44.977753
$ws.Range("F10").Value = -93.265011000000001
$ws.Range("E11").Value = 39.114052999999998
$ws.Range("F11").Value = -94.627464000000003
$ws.Range("E12").Value = 40.712783999999999
$ws.Range("F12").Value = -74.005941000000007
$ws.Range("E13").Value = 32.776663999999997
$ws.Range("F13").Value = -96.796987999999999
$ws.Range("E14").Value = 35.467559999999999
$ws.Range("F14").Value = -97.516428000000005
$ws.Range("E15").Value = 35.085334000000003
$ws.Range("F15").Value = -106.605553
$ws.Range("E16").Value = 30.267153
$ws.Range("F16").Value = -97.743060999999997
$ws.Range("E17").Value = 33.448377000000001
$ws.Range("F17").Value = -112.074037
$ws.Range("E18").Value = 36.169941000000001
$ws.Range("F18").Value = -115.13983
$ws.Range("E19").Value = 34.052233999999999
$ws.Range("F19").Value = -118.243685
$ws.Range("E20").Value = 37.774929
$ws.Range("F20").Value = -122.419416
$ws.Range("E21").Value = 45.523062000000003
$ws.Range("F21").Value = -122.67648199999999
$ws.Range("E22").Value = 47.606209
$ws.Range("F22").Value = -122.332071
$ws.Range("E23").Value = 42.280825999999998
$ws.Range("F23").Value = -83.743037999999999
$ws.Range("E24").Value = 39.952584000000002
$ws.Range("F24").Value = -75.165222
$ws.Range("E25").Value = 41.823988999999997
$ws.Range("F25").Value = -71.412834000000004
$ws.Range("E26").Value = 29.951066000000001
$ws.Range("F26").Value = -90.071532000000005
$ws.Range("E27").Value = 25.761679999999998
$ws.Range("F27").Value = -80.191789999999997
$ws.Range("E28").Value = 29.424122000000001
$ws.Range("F28").Value = -98.493628000000001
$ws.Range("E29").Value = 37.687176000000001
$ws.Range("F29").Value = -97.330053000000007
$ws.Range("E30").Value = 39.768402999999999
$ws.Range("F30").Value = -86.158068
$ws.Range("E31").Value = 27.950575000000001
$ws.Range("F31").Value = -82.457177999999999

# Update selection to match the target state (F18 was last active cell)
[void]$ws.Range("F18").Select()
